$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value2 = 13.22947826302449
$ws.Cells.Item(2, 3).Value2 = 11.05876974005799
$ws.Cells.Item(2, 4).Value2 = 9.57220978971997
$ws.Cells.Item(2, 6).Value2 = 34.26146669669116
$ws.Cells.Item(2, 7).Value2 = 36.07939500997963
$ws.Cells.Item(2, 8).Value2 = 16.28047994266157
$ws.Cells.Item(2, 10).Value2 = 10.90908626391087
$ws.Cells.Item(2, 11).Value2 = 9.170600578652673
$ws.Cells.Item(2, 12).Value2 = 11.2279975162594
$ws.Cells.Item(2, 15).Value2 = 25.72290425707454
$ws.Cells.Item(3, 2).Value2 = 12.97435284918561
$ws.Cells.Item(3, 3).Value2 = 11.08100097019673
$ws.Cells.Item(3, 4).Value2 = 9.551659272247031
$ws.Cells.Item(3, 6).Value2 = 34.36554213281768
$ws.Cells.Item(3, 7).Value2 = 36.22348191828203
$ws.Cells.Item(3, 8).Value2 = 16.33406521253697
$ws.Cells.Item(3, 10).Value2 = 10.93454643410926
$ws.Cells.Item(3, 11).Value2 = 8.971254705129171
$ws.Cells.Item(3, 12).Value2 = 11.21886028932719
$ws.Cells.Item(3, 15).Value2 = 25.81909088076728
$ws.Cells.Item(4, 2).Value2 = 12.81683077433349
$ws.Cells.Item(4, 3).Value2 = 11.09560576612598
$ws.Cells.Item(4, 4).Value2 = 9.540397122944148
$ws.Cells.Item(4, 6).Value2 = 34.43631281885576
$ws.Cells.Item(4, 7).Value2 = 36.32102793603581
$ws.Cells.Item(4, 8).Value2 = 16.36921147705213
$ws.Cells.Item(4, 10).Value2 = 10.95118637905992
$ws.Cells.Item(4, 11).Value2 = 8.847357716119959
$ws.Cells.Item(4, 12).Value2 = 11.21459648685309
$ws.Cells.Item(4, 15).Value2 = 25.88274567337351
$ws.Cells.Item(5, 2).Value2 = 12.75250255047548
$ws.Cells.Item(5, 3).Value2 = 11.10179796306167
$ws.Cells.Item(5, 4).Value2 = 9.536152100853606
$ws.Cells.Item(5, 6).Value2 = 34.46687739018469
$ws.Cells.Item(5, 7).Value2 = 36.36305495342317
$ws.Cells.Item(5, 8).Value2 = 16.38409884644246
$ws.Cells.Item(5, 10).Value2 = 10.95822110455121
$ws.Cells.Item(5, 11).Value2 = 8.796557141934267
$ws.Cells.Item(5, 12).Value2 = 11.21319946167113
$ws.Cells.Item(5, 15).Value2 = 25.90984077040241
$ws.Cells.Item(6, 2).Value2 = 12.74181510928964
$ws.Cells.Item(6, 3).Value2 = 11.10284072195226
$ws.Cells.Item(6, 4).Value2 = 9.535468114238519
$ws.Cells.Item(6, 6).Value2 = 34.47205673341347
$ws.Cells.Item(6, 7).Value2 = 36.37017080645229
$ws.Cells.Item(6, 8).Value2 = 16.38660502217787
$ws.Cells.Item(6, 10).Value2 = 10.95940456035173
$ws.Cells.Item(6, 11).Value2 = 8.78810488025303
$ws.Cells.Item(6, 12).Value2 = 11.21298810610805
$ws.Cells.Item(6, 15).Value2 = 25.91440965340756
$ws.Cells.Item(7, 2).Value2 = 12.81596366362601
$ws.Cells.Item(7, 3).Value2 = 11.09568830124733
$ws.Cells.Item(7, 4).Value2 = 9.540338474291502
$ws.Cells.Item(7, 6).Value2 = 34.43671804194903
$ws.Cells.Item(7, 7).Value2 = 36.32158551884449
$ws.Cells.Item(7, 8).Value2 = 16.36940996501898
$ws.Cells.Item(7, 10).Value2 = 10.95128022345564
$ws.Cells.Item(7, 11).Value2 = 8.846673775385394
$ws.Cells.Item(7, 12).Value2 = 11.21457626493259
$ws.Cells.Item(7, 15).Value2 = 25.88310641028471
$ws.Cells.Item(8, 2).Value2 = 13.14174160168423
$ws.Cells.Item(8, 3).Value2 = 11.06623724331894
$ws.Cells.Item(8, 4).Value2 = 9.564844671620095
$ws.Cells.Item(8, 6).Value2 = 34.29592526761599
$ws.Cells.Item(8, 7).Value2 = 36.12718886423164
$ws.Cells.Item(8, 8).Value2 = 16.29849057965371
$ws.Cells.Item(8, 10).Value2 = 10.91765617872283
$ws.Cells.Item(8, 11).Value2 = 9.102215236650192
$ws.Cells.Item(8, 12).Value2 = 11.22456872713371
$ws.Cells.Item(8, 15).Value2 = 25.75511530323302
$ws.Cells.Item(9, 2).Value2 = 13.77012582586411
$ws.Cells.Item(9, 3).Value2 = 11.016034724283
$ws.Cells.Item(9, 4).Value2 = 9.623499875545628
$ws.Cells.Item(9, 6).Value2 = 34.07441266496973
$ws.Cells.Item(9, 7).Value2 = 35.81825484290141
$ws.Cells.Item(9, 8).Value2 = 16.17720430074247
$ws.Cells.Item(9, 10).Value2 = 10.85968852795131
$ws.Cells.Item(9, 11).Value2 = 9.588726116058551
$ws.Cells.Item(9, 12).Value2 = 11.25475886102277
$ws.Cells.Item(9, 15).Value2 = 25.54060857860528
$ws.Cells.Item(10, 2).Value2 = 14.22089750755456
$ws.Cells.Item(10, 3).Value2 = 10.9837208351198
$ws.Cells.Item(10, 4).Value2 = 9.672826677338191
$ws.Cells.Item(10, 6).Value2 = 33.94505444143533
$ws.Cells.Item(10, 7).Value2 = 35.63569168833312
$ws.Cells.Item(10, 8).Value2 = 16.09890478974659
$ws.Cells.Item(10, 10).Value2 = 10.82192535076237
$ws.Cells.Item(10, 11).Value2 = 9.933836025789601
$ws.Cells.Item(10, 12).Value2 = 11.28327289086737
$ws.Cells.Item(10, 15).Value2 = 25.40527549294825
$ws.Cells.Item(11, 2).Value2 = 14.42270894667345
$ws.Cells.Item(11, 3).Value2 = 10.97000572746945
$ws.Cells.Item(11, 4).Value2 = 9.696566729611039
$ws.Cells.Item(11, 6).Value2 = 33.89347676555048
$ws.Cells.Item(11, 7).Value2 = 35.56235289933957
$ws.Cells.Item(11, 8).Value2 = 16.06562458833897
$ws.Cells.Item(11, 10).Value2 = 10.80578690551399
$ws.Cells.Item(11, 11).Value2 = 10.08750040888381
$ws.Cells.Item(11, 12).Value2 = 11.29759033505864
$ws.Cells.Item(11, 15).Value2 = 25.34854891503797
$ws.Cells.Item(12, 2).Value2 = 14.49859388727008
$ws.Cells.Item(12, 3).Value2 = 10.96495322623205
$ws.Cells.Item(12, 4).Value2 = 9.705738669532249
$ws.Cells.Item(12, 6).Value2 = 33.87499224143102
$ws.Cells.Item(12, 7).Value2 = 35.53598319961246
$ws.Cells.Item(12, 8).Value2 = 16.05335806236167
$ws.Cells.Item(12, 10).Value2 = 10.79982476220932
$ws.Cells.Item(12, 11).Value2 = 10.14516051541222
$ws.Cells.Item(12, 12).Value2 = 11.30320280768227
$ws.Cells.Item(12, 15).Value2 = 25.32776415536865
$ws.Cells.Item(13, 2).Value2 = 14.48227564277014
$ws.Cells.Item(13, 3).Value2 = 10.96603510431551
$ws.Cells.Item(13, 4).Value2 = 9.703755311010326
$ws.Cells.Item(13, 6).Value2 = 33.87892663192633
$ws.Cells.Item(13, 7).Value2 = 35.54159993752389
$ws.Cells.Item(13, 8).Value2 = 16.05598494123032
$ws.Cells.Item(13, 10).Value2 = 10.80110219064238
$ws.Cells.Item(13, 11).Value2 = 10.13276668755898
$ws.Cells.Item(13, 12).Value2 = 11.30198562330438
$ws.Cells.Item(13, 15).Value2 = 25.33220954163888
$ws.Cells.Item(14, 2).Value2 = 14.42896318694271
$ws.Cells.Item(14, 3).Value2 = 10.96958722976745
$ws.Cells.Item(14, 4).Value2 = 9.697317690451454
$ws.Cells.Item(14, 6).Value2 = 33.8919350437706
$ws.Cells.Item(14, 7).Value2 = 35.56015531562849
$ws.Cells.Item(14, 8).Value2 = 16.06460868248755
$ws.Cells.Item(14, 10).Value2 = 10.80529341002996
$ws.Cells.Item(14, 11).Value2 = 10.09225502525444
$ws.Cells.Item(14, 12).Value2 = 11.29804826996932
$ws.Cells.Item(14, 15).Value2 = 25.34682497963231
$ws.Cells.Item(15, 2).Value2 = 14.39623587817808
$ws.Cells.Item(15, 3).Value2 = 10.97178137207889
$ws.Cells.Item(15, 4).Value2 = 9.693398024750472
$ws.Cells.Item(15, 6).Value2 = 33.90003944559006
$ws.Cells.Item(15, 7).Value2 = 35.5717037841357
$ws.Cells.Item(15, 8).Value2 = 16.06993471588244
$ws.Cells.Item(15, 10).Value2 = 10.80788006164779
$ws.Cells.Item(15, 11).Value2 = 10.06737009405119
$ws.Cells.Item(15, 12).Value2 = 11.29566128469851
$ws.Cells.Item(15, 15).Value2 = 25.35586806812929
$ws.Cells.Item(16, 2).Value2 = 14.20763726076375
$ws.Cells.Item(16, 3).Value2 = 10.98463692219143
$ws.Cells.Item(16, 4).Value2 = 9.671300967547806
$ws.Cells.Item(16, 6).Value2 = 33.94857156396166
$ws.Cells.Item(16, 7).Value2 = 35.64068044520987
$ws.Cells.Item(16, 8).Value2 = 16.1011267601256
$ws.Cells.Item(16, 10).Value2 = 10.82300092964451
$ws.Cells.Item(16, 11).Value2 = 9.923722313128405
$ws.Cells.Item(16, 12).Value2 = 11.28236406430927
$ws.Cells.Item(16, 15).Value2 = 25.409080111146
$ws.Cells.Item(17, 2).Value2 = 14.09105555988985
$ws.Cells.Item(17, 3).Value2 = 10.99277523442534
$ws.Cells.Item(17, 4).Value2 = 9.658074882976738
$ws.Cells.Item(17, 6).Value2 = 33.98020722600747
$ws.Cells.Item(17, 7).Value2 = 35.68548659817129
$ws.Cells.Item(17, 8).Value2 = 16.12086081213509
$ws.Cells.Item(17, 10).Value2 = 10.83254319720426
$ws.Cells.Item(17, 11).Value2 = 9.834709464921648
$ws.Cells.Item(17, 12).Value2 = 11.27454949660424
$ws.Cells.Item(17, 15).Value2 = 25.44296344185995
$ws.Cells.Item(18, 2).Value2 = 14.02369915344597
$ws.Cells.Item(18, 3).Value2 = 10.99754888381192
$ws.Cells.Item(18, 4).Value2 = 9.65059036161496
$ws.Cells.Item(18, 6).Value2 = 33.99908723439998
$ws.Cells.Item(18, 7).Value2 = 35.71217154224396
$ws.Cells.Item(18, 8).Value2 = 16.13243145954945
$ws.Cells.Item(18, 10).Value2 = 10.83812959326194
$ws.Cells.Item(18, 11).Value2 = 9.783201197659169
$ws.Cells.Item(18, 12).Value2 = 11.27018165906824
$ws.Cells.Item(18, 15).Value2 = 25.46290748273959
$ws.Cells.Item(19, 2).Value2 = 14.00084380904839
$ws.Cells.Item(19, 3).Value2 = 10.9991810967229
$ws.Cells.Item(19, 4).Value2 = 9.648077466221855
$ws.Cells.Item(19, 6).Value2 = 34.00559710600157
$ws.Cells.Item(19, 7).Value2 = 35.72136334142157
$ws.Cells.Item(19, 8).Value2 = 16.13638690191312
$ws.Cells.Item(19, 10).Value2 = 10.84003788527369
$ws.Cells.Item(19, 11).Value2 = 9.765709671793154
$ws.Cells.Item(19, 12).Value2 = 11.26872466247021
$ws.Cells.Item(19, 15).Value2 = 25.46973835701156
$ws.Cells.Item(20, 2).Value2 = 14.10349762396238
$ws.Cells.Item(20, 3).Value2 = 10.99189930580511
$ws.Cells.Item(20, 4).Value2 = 9.659470151242427
$ws.Cells.Item(20, 6).Value2 = 33.97676875021601
$ws.Cells.Item(20, 7).Value2 = 35.68062230651681
$ws.Cells.Item(20, 8).Value2 = 16.11873730709396
$ws.Cells.Item(20, 10).Value2 = 10.83151727375582
$ws.Cells.Item(20, 11).Value2 = 9.844217542028838
$ws.Cells.Item(20, 12).Value2 = 11.27536825753242
$ws.Cells.Item(20, 15).Value2 = 25.43930938154974
$ws.Cells.Item(21, 2).Value2 = 14.44463744803731
$ws.Cells.Item(21, 3).Value2 = 10.96854005844977
$ws.Cells.Item(21, 4).Value2 = 9.699203674394973
$ws.Cells.Item(21, 6).Value2 = 33.88808573493868
$ws.Cells.Item(21, 7).Value2 = 35.55466705352958
$ws.Cells.Item(21, 8).Value2 = 16.06206656523401
$ws.Cells.Item(21, 10).Value2 = 10.80405830332738
$ws.Cells.Item(21, 11).Value2 = 10.10416903257192
$ws.Cells.Item(21, 12).Value2 = 11.29919961191432
$ws.Cells.Item(21, 15).Value2 = 25.34251316576547
$ws.Cells.Item(22, 2).Value2 = 14.66443450919456
$ws.Cells.Item(22, 3).Value2 = 10.95409573306573
$ws.Cells.Item(22, 4).Value2 = 9.726230772411197
$ws.Cells.Item(22, 6).Value2 = 33.83622858318914
$ws.Cells.Item(22, 7).Value2 = 35.48052307489782
$ws.Cells.Item(22, 8).Value2 = 16.02698707079787
$ws.Cells.Item(22, 10).Value2 = 10.78698135236271
$ws.Cells.Item(22, 11).Value2 = 10.27095646956198
$ws.Cells.Item(22, 12).Value2 = 11.31588531004004
$ws.Cells.Item(22, 15).Value2 = 25.28331055512667
$ws.Cells.Item(23, 2).Value2 = 14.54743556512819
$ws.Cells.Item(23, 3).Value2 = 10.96172985795011
$ws.Cells.Item(23, 4).Value2 = 9.711710709993515
$ws.Cells.Item(23, 6).Value2 = 33.86334685092128
$ws.Cells.Item(23, 7).Value2 = 35.51934524448081
$ws.Cells.Item(23, 8).Value2 = 16.04553059546099
$ws.Cells.Item(23, 10).Value2 = 10.79601626813729
$ws.Cells.Item(23, 11).Value2 = 10.18223887831515
$ws.Cells.Item(23, 12).Value2 = 11.30687919399426
$ws.Cells.Item(23, 15).Value2 = 25.31453641333204
$ws.Cells.Item(24, 2).Value2 = 14.09787359565825
$ws.Cells.Item(24, 3).Value2 = 10.99229501772733
$ws.Cells.Item(24, 4).Value2 = 9.65883897839195
$ws.Cells.Item(24, 6).Value2 = 33.97832112896423
$ws.Cells.Item(24, 7).Value2 = 35.68281857100462
$ws.Cells.Item(24, 8).Value2 = 16.11969664196815
$ws.Cells.Item(24, 10).Value2 = 10.83198078091697
$ws.Cells.Item(24, 11).Value2 = 9.839919975553183
$ws.Cells.Item(24, 12).Value2 = 11.27499770658056
$ws.Cells.Item(24, 15).Value2 = 25.44095993674792
$ws.Cells.Item(25, 2).Value2 = 13.60173142717086
$ws.Cells.Item(25, 3).Value2 = 11.02881092694039
$ws.Cells.Item(25, 4).Value2 = 9.60652006598008
$ws.Cells.Item(25, 6).Value2 = 34.12848208150153
$ws.Cells.Item(25, 7).Value2 = 35.89405510508868
$ws.Cells.Item(25, 8).Value2 = 16.2081148120579
$ws.Cells.Item(25, 10).Value2 = 10.87452055201091
$ws.Cells.Item(25, 11).Value2 = 9.459042193630546
$ws.Cells.Item(25, 12).Value2 = 11.24546976929402
$ws.Cells.Item(25, 15).Value2 = 25.5947301044241
